$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Drop the "Year 2 Candidates attending the Catholic HS teens
#    sessions" banner row entirely - rows below it shift up by one.
# ------------------------------------------------------------------
$ws.Rows.Item(5).Delete()

# ------------------------------------------------------------------
# 2. Build the bold/underlined (no italic), bordered, filled header
#    look for D1 from the existing title cell (A1) before we touch it.
# ------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Font.Italic = $false

# ------------------------------------------------------------------
# 3. Build the bold/underlined, bordered, UNFILLED header look (for
#    A1, B1, C1, E1) from the plain bordered cell B1, then layer on
#    bold + underline + text format.
# ------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").Font.Underline = 2
$ws.Range("E1").NumberFormat = "@"

$ws.Range("B1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").Font.Underline = 2
$ws.Range("A1:C1").NumberFormat = "@"

# ------------------------------------------------------------------
# 4. New header row text + the new "Cardinal Gibbons HS Group" label.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Last Name"
$ws.Range("B1").Value = "1st name"
$ws.Range("C1").Value = "grade "
$ws.Range("D1").Value = "parents preferred email"
$ws.Range("E1").Value = "Cardinal Gibbons HS Group"

# ------------------------------------------------------------------
# 5. Typo fix: "@nc.rr.com; rannunz" -> "@nc.rr.com, rannunz"
#    (this is row 5 now, after the row-5 delete shifted rows up).
# ------------------------------------------------------------------
$ws.Range("D5").Value = "@nc.rr.com, rannunz"

# ------------------------------------------------------------------
# 6. Stamp "Cardinal Gibbons HS Group" on E5/E6, matching the plain
#    bordered/unfilled look used elsewhere in the data rows, with a
#    text number format.
# ------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("E5:E6").PasteSpecial(-4122)
$ws.Range("E5:E6").NumberFormat = "@"
$ws.Range("E5").Value = "Cardinal Gibbons HS Group"
$ws.Range("E6").Value = "Cardinal Gibbons HS Group"
